$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
  # Row 11
$ws.Range("H11").Value = 90918.17999999999
$ws.Range("I11").Value = 90918.17999999999
$ws.Range("K11").Value = 90918.17999999999
$ws.Range("M11").Value = -90778.17999999999
  # Row 46
$ws.Range("H46").Value = 958.1667
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 958.1667
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 2874.5001
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -3112.5001
  # Row 60
$ws.Range("H60").Value = 958.1667
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 958.1667
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 2874.5001
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = -3842.5001
  # Row 64
$ws.Range("H64").Value = 69331.8
$ws.Range("I64").Value = 127437.375
$ws.Range("J64").Value = 2925.4285
$ws.Range("K64").Value = 127437.375
$ws.Range("L64").Value = 2925.4285
$ws.Range("M64").Value = -127189.375
$ws.Range("N64").Value = -3421.4285
  # Row 67
$ws.Range("H67").Value = 69331.8
$ws.Range("I67").Value = 127437.375
$ws.Range("J67").Value = 2925.4285
$ws.Range("K67").Value = 127437.375
$ws.Range("L67").Value = 2925.4285
$ws.Range("M67").Value = -126579.375
$ws.Range("N67").Value = -4641.4285
  # Row 70
$ws.Range("H70").Value = 1140.4333
$ws.Range("I70").Value = 1258.8096
$ws.Range("J70").Value = 864.2222
$ws.Range("K70").Value = 3776.4288
$ws.Range("L70").Value = 2592.6666
$ws.Range("M70").Value = -3506.4288
$ws.Range("N70").Value = -3132.6666
  # Row 73
$ws.Range("H73").Value = 1140.4333
$ws.Range("I73").Value = 1258.8096
$ws.Range("J73").Value = 864.2222
$ws.Range("K73").Value = 3776.4288
$ws.Range("L73").Value = 2592.6666
$ws.Range("M73").Value = -2840.4288
$ws.Range("N73").Value = -4464.6666
  # Row 129
$ws.Range("H129").Value = 2619.4138
$ws.Range("I129").Value = 8230.691999999999
$ws.Range("J129").Value = 998.3778
$ws.Range("K129").Value = 24692.076
$ws.Range("L129").Value = 2995.1334
$ws.Range("M129").Value = -19692.076
$ws.Range("N129").Value = -12995.1334
  # Row 132
$ws.Range("H132").Value = 8628155
$ws.Range("I132").Value = 9623596
$ws.Range("K132").Value = 28870788
$ws.Range("M132").Value = -28868258
  # Row 138
$ws.Range("H138").Value = 3914.6
$ws.Range("I138").Value = 3424.25
$ws.Range("J138").Value = 4145.353
$ws.Range("K138").Value = 10272.75
$ws.Range("L138").Value = 12436.059
$ws.Range("M138").Value = -5132.75
$ws.Range("N138").Value = -22716.059

$ws = $wb.Worksheets.Item("ARM")
  # Row 36
$ws.Range("H36").Value = 700
$ws.Range("I36").Value = 700
$ws.Range("K36").Value = 700
$ws.Range("M36").Value = -354
  # Row 102
$ws.Range("H102").Value = 73595.36
$ws.Range("I102").Value = 126923.75
$ws.Range("J102").Value = 2490.8333
$ws.Range("K102").Value = 126923.75
$ws.Range("L102").Value = 2490.8333
$ws.Range("M102").Value = -125301.75
$ws.Range("N102").Value = -5734.8333
  # Row 122
$ws.Range("H122").Value = 1186.0834
$ws.Range("I122").Value = 1189.875
$ws.Range("J122").Value = 1178.5
$ws.Range("K122").Value = 3569.625
$ws.Range("L122").Value = 3535.5
$ws.Range("M122").Value = -1119.625
$ws.Range("N122").Value = -8435.5

$ws = $wb.Worksheets.Item("BSM")
  # Row 43
$ws.Range("H43").Value = 398000
$ws.Range("J43").Value = 398000
$ws.Range("L43").Value = 398000
$ws.Range("N43").Value = -398362

$ws = $wb.Worksheets.Item("CRP")
  # Row 7
$ws.Range("H7").Value = 139.6875
$ws.Range("I7").Value = 101.454544
$ws.Range("J7").Value = 223.8
$ws.Range("K7").Value = 101.454544
$ws.Range("L7").Value = 223.8
$ws.Range("M7").Value = 11.545456
$ws.Range("N7").Value = -449.8
  # Row 51
$ws.Range("H51").Value = 7983
$ws.Range("J51").Value = 7983
$ws.Range("L51").Value = 7983
$ws.Range("N51").Value = -9455
  # Row 58
$ws.Range("H58").Value = 1590.8422
$ws.Range("I58").Value = 1491.5834
$ws.Range("J58").Value = 1761
$ws.Range("K58").Value = 1491.5834
$ws.Range("L58").Value = 1761
$ws.Range("M58").Value = -1288.5834
$ws.Range("N58").Value = -2167
  # Row 61
$ws.Range("H61").Value = 7983
$ws.Range("J61").Value = 7983
$ws.Range("L61").Value = 7983
$ws.Range("N61").Value = -8679
  # Row 99
$ws.Range("H99").Value = 45757
$ws.Range("I99").Value = 3080
$ws.Range("J99").Value = 57950.43
$ws.Range("K99").Value = 3080
$ws.Range("L99").Value = 57950.43
$ws.Range("M99").Value = -1582
$ws.Range("N99").Value = -60946.43
  # Row 126
$ws.Range("H126").Value = 45757
$ws.Range("I126").Value = 3080
$ws.Range("J126").Value = 57950.43
$ws.Range("K126").Value = 9240
$ws.Range("L126").Value = 173851.29
$ws.Range("M126").Value = -6770
$ws.Range("N126").Value = -178791.29
  # Row 136
$ws.Range("H136").Value = 1590.8422
$ws.Range("I136").Value = 1491.5834
$ws.Range("J136").Value = 1761
$ws.Range("K136").Value = 4474.7502
$ws.Range("L136").Value = 5283
$ws.Range("M136").Value = -1924.7502
$ws.Range("N136").Value = -10383

$ws = $wb.Worksheets.Item("CUL")
  # Row 49
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()
  # Row 117
$ws.Range("H117").Value = 6494.737
$ws.Range("J117").Value = 7157
$ws.Range("L117").Value = 21471
$ws.Range("N117").Value = -28355
  # Row 131
$ws.Range("H131").Value = 798.3099999999999
$ws.Range("J131").Value = 811.3196
$ws.Range("L131").Value = 2433.9588
$ws.Range("N131").Value = -12513.9588
  # Row 138
$ws.Range("H138").Value = 12292.9
$ws.Range("I138").Value = 16275.714
$ws.Range("K138").Value = 48827.142
$ws.Range("M138").Value = -43687.142

$ws = $wb.Worksheets.Item("GSM")
  # Row 113
$ws.Range("H113").Value = 1482.3077
$ws.Range("I113").Value = 833
$ws.Range("J113").Value = 1600.3636
$ws.Range("K113").Value = 833
$ws.Range("L113").Value = 1600.3636
$ws.Range("M113").Value = 1337
$ws.Range("N113").Value = -5940.3636
  # Row 132
$ws.Range("H132").Value = 3081.0952
$ws.Range("J132").Value = 3250.1428
$ws.Range("L132").Value = 9750.428400000001
$ws.Range("N132").Value = -14810.4284

$ws = $wb.Worksheets.Item("LTW")
  # Row 7
$ws.Range("H7").Value = 6534.727
$ws.Range("I7").Value = 7834
$ws.Range("J7").Value = 4975.6
$ws.Range("K7").Value = 7834
$ws.Range("L7").Value = 4975.6
$ws.Range("M7").Value = -7722
$ws.Range("N7").Value = -5199.6
  # Row 93
$ws.Range("H93").Value = 2089.95
$ws.Range("I93").Value = 2819.3
$ws.Range("J93").Value = 1360.6
$ws.Range("K93").Value = 2819.3
$ws.Range("L93").Value = 1360.6
$ws.Range("M93").Value = -1571.3
$ws.Range("N93").Value = -3856.6
  # Row 126
$ws.Range("H126").Value = 6534.727
$ws.Range("I126").Value = 7834
$ws.Range("J126").Value = 4975.6
$ws.Range("K126").Value = 23502
$ws.Range("L126").Value = 14926.8
$ws.Range("M126").Value = -21032
$ws.Range("N126").Value = -19866.8
